{"js": "const doc = context.document;\nconst body = doc.body;\n\n// 1. Remove the existing \"_GoBack\" bookmark (Word keeps only one such\n//    bookmark, marking the location of the last edit; it will be re-created\n//    below at the new edit location).\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Replace the old folder name with the new one.\nconst hits = body.search(\"day2-files\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"FrontEnd/Code/Day 2\", \"Replace\");\n  await context.sync();\n}\n\n// 3. Drop a new \"_GoBack\" bookmark right after \"FrontEnd/Code/D\" -- the spot\n//    Word left its last-edit mark at.\nconst anchorHits = body.search(\"FrontEnd/Code/D\", { matchCase: true });\nanchorHits.load(\"items\");\nawait context.sync();\n\nif (anchorHits.items.length > 0) {\n  const afterPrefix = anchorHits.items[0].getRange(\"End\");\n  afterPrefix.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the existing \"_GoBack\" bookmark (it will be re-created at the new\n#    edit location below, mirroring what Word itself does when the last edit\n#    moves elsewhere in the document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Find the old folder name and replace it with the new one.\n#    wdFindContinue (1), Forward, Replace:=wdReplaceOne (1).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"day2-files\", $false, $false, $false, $false, $false, $true, 1, $false, \"FrontEnd/Code/Day 2\", 1) | Out-Null\n\n# 3. Re-find the newly inserted text and drop a \"_GoBack\" bookmark right after\n#    \"FrontEnd/Code/D\" (the position Word left its last-edit mark at).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"FrontEnd/Code/D\"\n$find2.Execute() | Out-Null\n\n$mark = $find2.Parent.Duplicate\n$mark.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $mark)\n"}
